# NIT-9014182136.xlsx - "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The worker detail table (rows 16-73, columns C:G) is rebuilt: instead of
# being grouped period-by-period (2202..2208) with every worker repeated in
# each period block, it is regrouped worker-by-worker, each worker keeping
# their own contiguous block of periods ordered from the newest (2208) down
# to the oldest (2202). Each worker keeps the same "Valor Mora" (F) /
# "Salario Basico" (G) figures they already had; only their position in the
# table (and therefore which period lands on which row) changes. The last
# worker (JUAN CARLOS OCHOA PEREZ) only has two periods of mora (2204, 2203).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Doc Number, Worker name, then a list of (Periodo, ValorMora, SalarioBasico)
$workers = @(
    @{ Doc = "8980918";    Name = "JUAN BAUTISTA TIRADO FERIA";        Periods = @(
        @{ P = "2208"; F = 29333; G = 1000000 },
        @{ P = "2207"; F = 40000; G = 1000000 },
        @{ P = "2206"; F = 40000; G = 1000000 },
        @{ P = "2205"; F = 40000; G = 1000000 },
        @{ P = "2204"; F = 40000; G = 1000000 },
        @{ P = "2203"; F = 40000; G = 1000000 },
        @{ P = "2202"; F = 40000; G = 1000000 }
    ) },
    @{ Doc = "77164702";   Name = "JOSE LUIS DE HORTA VAZQUEZ";        Periods = @(
        @{ P = "2208"; F = 29333; G = 1000000 },
        @{ P = "2207"; F = 40000; G = 1000000 },
        @{ P = "2206"; F = 40000; G = 1000000 },
        @{ P = "2205"; F = 40000; G = 1000000 },
        @{ P = "2204"; F = 40000; G = 1000000 },
        @{ P = "2203"; F = 40000; G = 1000000 },
        @{ P = "2202"; F = 40000; G = 1000000 }
    ) },
    @{ Doc = "8850786";    Name = "ADOLFO ANGEL TIRADO FERIA";         Periods = @(
        @{ P = "2208"; F = 29333; G = 1000000 },
        @{ P = "2207"; F = 40000; G = 1000000 },
        @{ P = "2206"; F = 40000; G = 1000000 },
        @{ P = "2205"; F = 40000; G = 1000000 },
        @{ P = "2204"; F = 40000; G = 1000000 },
        @{ P = "2203"; F = 40000; G = 1000000 },
        @{ P = "2202"; F = 40000; G = 1000000 }
    ) },
    @{ Doc = "12684694";   Name = "LIBARDO DE JESUS VARGAS CARVAJAL";  Periods = @(
        @{ P = "2208"; F = 58667; G = 2000000 },
        @{ P = "2207"; F = 80000; G = 2000000 },
        @{ P = "2206"; F = 80000; G = 2000000 },
        @{ P = "2205"; F = 80000; G = 2000000 },
        @{ P = "2204"; F = 80000; G = 2000000 },
        @{ P = "2203"; F = 80000; G = 2000000 },
        @{ P = "2202"; F = 80000; G = 2000000 }
    ) },
    @{ Doc = "1022362050";  Name = "NATALIA NIÑO MIRANDA";             Periods = @(
        @{ P = "2208"; F = 29333; G = 1000000 },
        @{ P = "2207"; F = 40000; G = 1000000 },
        @{ P = "2206"; F = 40000; G = 1000000 },
        @{ P = "2205"; F = 40000; G = 1000000 },
        @{ P = "2204"; F = 40000; G = 1000000 },
        @{ P = "2203"; F = 40000; G = 1000000 },
        @{ P = "2202"; F = 40000; G = 1000000 }
    ) },
    @{ Doc = "1002046926";  Name = "CARLOS MANUEL MOLINA DE ARCO";     Periods = @(
        @{ P = "2208"; F = 29333; G = 1000000 },
        @{ P = "2207"; F = 40000; G = 1000000 },
        @{ P = "2206"; F = 40000; G = 1000000 },
        @{ P = "2205"; F = 40000; G = 1000000 },
        @{ P = "2204"; F = 40000; G = 1000000 },
        @{ P = "2203"; F = 40000; G = 1000000 },
        @{ P = "2202"; F = 40000; G = 1000000 }
    ) },
    @{ Doc = "1049795332";  Name = "YUBERT EDUARDO HERNANDEZ SOLER";   Periods = @(
        @{ P = "2208"; F = 29333; G = 1000000 },
        @{ P = "2207"; F = 40000; G = 1000000 },
        @{ P = "2206"; F = 40000; G = 1000000 },
        @{ P = "2205"; F = 40000; G = 1000000 },
        @{ P = "2204"; F = 40000; G = 1000000 },
        @{ P = "2203"; F = 40000; G = 1000000 },
        @{ P = "2202"; F = 40000; G = 1000000 }
    ) },
    @{ Doc = "92032673";    Name = "RAFAEL EDUARDO PINEDA CASTILLO";   Periods = @(
        @{ P = "2208"; F = 29333; G = 1000000 },
        @{ P = "2207"; F = 40000; G = 1000000 },
        @{ P = "2206"; F = 40000; G = 1000000 },
        @{ P = "2205"; F = 40000; G = 1000000 },
        @{ P = "2204"; F = 40000; G = 1000000 },
        @{ P = "2203"; F = 40000; G = 1000000 },
        @{ P = "2202"; F = 40000; G = 1000000 }
    ) },
    @{ Doc = "1042578744";  Name = "JUAN CARLOS OCHOA PEREZ";          Periods = @(
        @{ P = "2204"; F = 6667;  G = 1000000 },
        @{ P = "2203"; F = 40000; G = 1000000 }
    ) }
)

$row = 16
foreach ($worker in $workers) {
    foreach ($period in $worker.Periods) {
        $ws.Range("C$row").Value = $worker.Doc
        $ws.Range("D$row").Value = $worker.Name
        $ws.Range("E$row").Value = $period.P
        $ws.Range("F$row").Value = $period.F
        $ws.Range("G$row").Value = $period.G
        $row = $row + 1
    }
}
